$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the sample "TES-0987" book code with "TES-123" in the second sample row.
$ws.Range("D9").Value = "TES-123"

# Move the active selection to D8 (matches the saved sheet view state).
$ws.Range("D8").Select() | Out-Null
